$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Enter data for Day 4 (row 8): Part One (C), Part Two (D)
$ws.Range("C8").Value = 0.01037037037037037
$ws.Range("D8").Value = 0.008113425925925925
$ws.Range("E8").Formula = "=D8+C8"
$ws.Range("F8").Formula = "=E8/2"

# Set selection / view to match the after-state
$ws.Range("J22").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
